$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-12-22 02:02:09"
}
